$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '51.652.49'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '  +1.06%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '3.041.11'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '  +2.62%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '384.56'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +1.27%  '
$ws.Range('E6').Value = '  +0.61%  '
$ws.Range('E7').Value = '  -0.15%  '
$ws.Range('E8').Value = '  +0.03%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.589'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  -0.41%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '36.93'
$ws.Range('D10').Style = "Normal"
$ws.Range('E11').Value = '  +0.15%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.0863'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  +1.06%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '3.515.72'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  +2.58%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '18.70'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  +1.92%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '7.79'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  -0.44%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '3.053.43'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  +3.12%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.976'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  -2.45%  '
$ws.Range('E18').Value = '  -11.35%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '51.693.29'
$ws.Range('D19').Style = "Normal"
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '3.09'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -0.72%  '
$ws.Range('B21').Value = 'InternetComputer(DFINITY)'
$ws.Range('C21').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '12.37'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -0.55%  '
$ws.Range('B22').Value = 'ShibaInu'
$ws.Range('C22').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '0.0₃0964'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +0.32%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '69.93'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -0.33%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '267.01'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  -0.40%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '8.37'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  +6.00%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '7.45'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +3.48%  '
$ws.Range('E28').Value = '  +3.84%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '26.40'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  +1.89%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '0.108'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -2.58%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '10.31'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  -1.11%  '
$ws.Range('B33').Value = 'InjectiveProtocol'
$ws.Range('C33').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '34.18'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  -0.93%  '
$ws.Range('B34').Value = 'Toncoin'
$ws.Range('C34').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '2.07'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  +1.43%  '
$ws.Range('E35').Value = '  -1.01%  '
$ws.Range('E36').Value = '  +2.99%  '
$ws.Range('E37').Value = '  -0.13%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '3.36'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  +3.65%  '
$ws.Range('E39').Value = '  +5.01%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '17.00'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  +2.46%  '
$ws.Range('E41').Value = '  +1.48%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '128.19'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  +2.33%  '
$ws.Range('E43').Value = '  -0.39%  '
$ws.Range('E44').Value = '  +1.16%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '3.67'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +3.55%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '21.77'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  +0.86%  '
$ws.Range('B47').Value = 'ApeXProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '2.48'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  +4.41%  '
$ws.Range('B48').Value = 'WEMIXToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '2.10'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  +3.65%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '2.036.06'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -0.81%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '3.342.03'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  +2.69%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.210'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  +8.19%  '
